# Addressing problem of duplicate generators without a general "slice(1)"
# since data may be lost this way.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("generator_file")
$ws2 = $wb.Worksheets.Item("unit_file")
$ws3 = $wb.Worksheets.Item("plant_file")

# --- unit_file: add six new rows (plant 50489, generators C1/C2/C4/C5) ---
# Columns A & B first (plant_id / unit_id), row by row.
$ws2.Cells.Item(20,1).Value = 50489
$ws2.Cells.Item(20,2).Value = "C1"
$ws2.Cells.Item(21,1).Value = 50489
$ws2.Cells.Item(21,2).Value = "C2"
$ws2.Cells.Item(22,1).Value = 50489
$ws2.Cells.Item(22,2).Value = "C4"
$ws2.Cells.Item(23,1).Value = 50489
$ws2.Cells.Item(23,2).Value = "C5"
$ws2.Cells.Item(24,1).Value = 50489
$ws2.Cells.Item(24,2).Value = "C4"
$ws2.Cells.Item(25,1).Value = 50489
$ws2.Cells.Item(25,2).Value = "C5"

# Columns D & E next (column_to_update / update), row by row.
$ws2.Cells.Item(20,4).Value = "prime_mover"
$ws2.Cells.Item(20,5).Value = "GT"
$ws2.Cells.Item(21,4).Value = "prime_mover"
$ws2.Cells.Item(21,5).Value = "GT"
$ws2.Cells.Item(22,4).Value = "prime_mover"
$ws2.Cells.Item(22,5).Value = "GT"
$ws2.Cells.Item(23,4).Value = "prime_mover"
$ws2.Cells.Item(23,5).Value = "GT"
$ws2.Cells.Item(24,4).Value = "prop"
$ws2.Cells.Item(24,5).Value = 0
$ws2.Cells.Item(25,4).Value = "prop"
$ws2.Cells.Item(25,5).Value = 0

# --- generator_file: add a new row for plant 56032 ---
$ws1.Cells.Item(6,1).Value = 56032
$ws1.Cells.Item(6,3).Value = "keep_leading_zeroes"

# --- plant_file: drop the duplicate-prone primary_fuel_type / primary_fuel_category
#     manual corrections (rows 10-13) - this data can be lost with a blind slice(1) ---
$ws3.Rows("10:13").Delete()

# --- restore per-sheet selections (sheet1 stays the active tab) ---
$ws3.Range("B4").Select()
$ws2.Range("E25").Select()
$ws1.Range("C7").Select()

Write-Host "Edits applied"
